$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Preserve G:L fill-style formatting for the new rows before writing text ---
# Row 14 reuses row 8s "disappear" G/J cells (style 3 / style 2 fills); row 8 itself
# becomes blank (still keeps its fill) once the clue moves down to the new Book row.
$ws.Range("G8:L8").Copy($ws.Range("G14:L14"))
$ws.Range("G8").ClearContents()
$ws.Range("J8").ClearContents()

# Rows 9-13 are brand new; copy the blank-but-styled G:L block from row 7.
$ws.Range("G7:L7").Copy($ws.Range("G9:L9"))
$ws.Range("G7:L7").Copy($ws.Range("G10:L10"))
$ws.Range("G7:L7").Copy($ws.Range("G11:L11"))
$ws.Range("G7:L7").Copy($ws.Range("G12:L12"))
$ws.Range("G7:L7").Copy($ws.Range("G13:L13"))

# --- Cell value writes (Investigate Panel2 dialogue + clue rows) ---
$ws.Range("A2").Value = "He"
$ws.Range("B2").Value = "Sir, allow me to show you the Lord’s study."
$ws.Range("C2").Value = "He-Regular1"
$ws.Range("E2").Value = "Study"
$ws.Range("F2").Value = "Suspicious"
$ws.Range("L2").Value = "He-Sad"
$ws.Range("A3").Value = "Dee"
$ws.Range("B3").Value = "Much appreciated."
$ws.Range("C3").Value = "He-Regular1"
$ws.Range("E3").Value = "Study"
$ws.Range("B4").Value = "Steward He skillfully unlocked the door to the rear study."
$ws.Range("C4").Value = "He-Regular1"
$ws.Range("E4").Value = "Study"
$ws.Range("A5").Value = "He"
$ws.Range("B5").Value = "Please, both of you, come in."
$ws.Range("C5").Value = "He-Regular1"
$ws.Range("E5").Value = "Study"
$ws.Range("A6").Value = "He"
$ws.Range("B6").Value = "No one has entered this study since the Lord passed away."
$ws.Range("C6").Value = "He-Regular1"
$ws.Range("E6").Value = "Study"
$ws.Range("A7").Value = "He"
$ws.Range("B7").Value = "Everything inside should be just as he left it."
$ws.Range("C7").Value = "He-Regular1"
$ws.Range("D7").Value = "DialogueVocal"
$ws.Range("E7").Value = "Study"
$ws.Range("A8").Value = "Dee"
$ws.Range("B8").Value = "Did the Lord often spend time in this study?"
$ws.Range("C8").Value = "He-Regular1"
$ws.Range("D8").Value = "DialogueVocal"
$ws.Range("E8").Value = "Study"
$ws.Range("A9").Value = "He"
$ws.Range("B9").Value = "Yes. Aside from resting in his room or receiving guests in the main hall, he spent most of his time here—reading and drinking tea."
$ws.Range("C9").Value = "He-Regular1"
$ws.Range("D9").Value = "DialogueVocal"
$ws.Range("E9").Value = "Study"
$ws.Range("A10").Value = "Dee"
$ws.Range("B10").Value = "Very well. Let’s begin the investigation."
$ws.Range("C10").Value = "He-Regular1"
$ws.Range("D10").Value = "DialogueVocal"
$ws.Range("E10").Value = "Study"
$ws.Range("B11").Value = "Click on any area you find suspicious to gather clues."
$ws.Range("C11").Value = "He-Regular1"
$ws.Range("D11").Value = "DialogueVocal"
$ws.Range("E11").Value = "Study"
$ws.Range("A12").Value = "Investigate"
$ws.Range("B12").Value = "Desk"
$ws.Range("C12").Value = "Desk"
$ws.Range("D12").Value = "DialogueVocal"
$ws.Range("E12").Value = "Study"
$ws.Range("B13").Value = "Paper"
$ws.Range("C13").Value = "Paper"
$ws.Range("D13").Value = "DialogueVocal"
$ws.Range("E13").Value = "Study"
$ws.Range("B14").Value = "Book"
$ws.Range("C14").Value = "Book"
$ws.Range("D14").Value = "DialogueVocal"
$ws.Range("E14").Value = "Study"
$ws.Range("B15").Value = "End Investigation"
$ws.Range("C15").Value = "StoryScript13"
$ws.Range("D15").Value = "DialogueVocal"
$ws.Range("E15").Value = "Study"

# G14/J14 keep the moved "disappear" marker
$ws.Range("G14").Value = "disappear"
$ws.Range("J14").Value = "disappear"

# --- Row heights to match the wrapped dialogue text ---
$ws.Rows.Item(2).RowHeight = 17
$ws.Rows.Item(3).RowHeight = 17
$ws.Rows.Item(4).RowHeight = 17
$ws.Rows.Item(5).RowHeight = 17
$ws.Rows.Item(6).RowHeight = 34
$ws.Rows.Item(7).RowHeight = 17
$ws.Rows.Item(8).RowHeight = 17
$ws.Rows.Item(9).RowHeight = 51
$ws.Rows.Item(10).RowHeight = 17
$ws.Rows.Item(11).RowHeight = 17
$ws.Rows.Item(12).RowHeight = 17
$ws.Rows.Item(13).RowHeight = 17
$ws.Rows.Item(14).RowHeight = 17
$ws.Rows.Item(15).RowHeight = 17

# --- Final selection, matching the saved workbook state ---
$ws.Range("B20").Select()
